$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New employee data (row 4): Miroslav Havelka
$ws.Range("E4").Value = "Miroslav"
$ws.Range("F4").Value = "Havelka"
$ws.Range("D4").Value = "Ne"

# Full name built from given name + surname (mirrors A3's formula)
$ws.Range("A4").Formula = '=CONCAT(E4, " ", F4)'

# Generated login / e-mail (mirrors B3's formula)
$ws.Range("B4").Formula = '=CONCATENATE(LOWER(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE($E4,"á","a"),"č","c"),"ď","d"),"ě","e"),"é","e"),"í","i"),"ľ","l"),"ň","n"),"ó","o"),"ö","o"),"ř","r"),"š","s"),"ť","t"),"ú","u"),"ů","u"),"ý","y"),"ž","z"),"á","a"),"č","c"),"ď","d"),"ě","e"),"é","e"),"í","i"),"ľ","l"),"ň","n"),"ó","o"),"ö","o"),"ř","r"),"š","s"),"ť","t"),"ú","u"),"ů","u"),"ý","y"),"ž","z")), "_", LOWER(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE($F4,"á","a"),"č","c"),"ď","d"),"ě","e"),"é","e"),"í","i"),"ľ","l"),"ň","n"),"ó","o"),"ö","o"),"ř","r"),"š","s"),"ť","t"),"ú","u"),"ů","u"),"ý","y"),"ž","z"),"á","a"),"č","c"),"ď","d"),"ě","e"),"é","e"),"í","i"),"ľ","l"),"ň","n"),"ó","o"),"ö","o"),"ř","r"),"š","s"),"ť","t"),"ú","u"),"ů","u"),"ý","y"),"ž","z")), "@mblab.cloud")'

# Random temporary-password-like token (mirrors C3's formula)
$ws.Range("C4").Formula = '=CHAR(RANDBETWEEN(65,90))&CHAR(42)&CHAR(RANDBETWEEN(65,90))&RANDBETWEEN(10,99)&RANDBETWEEN(10,99)&CHAR(35)&CHAR(RANDBETWEEN(65,90))&RANDBETWEEN(10,99)'

# Clear out the stray phone-number value that was left in P3
$ws.Range("P3").ClearContents()

# Move the active selection, as last left by the editing user
$ws.Range("C9").Select() | Out-Null
